# Trade #60 closed at 2026-02-17 12:52:56 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.1
$summary.Range("B4").Value = 0.09
$summary.Range("B5").Value = 0.03
$summary.Range("B6").Value = 60
$summary.Range("B7").Value = 26
$summary.Range("B9").Value = 43.33

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.1
$status.Range("D4").Value = 60
$status.Range("E4").Value = 0.09
$status.Range("F4").Value = 0.1
$status.Range("G4").Value = 43.33

# ---- New closed trade row, appended to both "All Trades" and "MarketMaking" ----
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A61").Value = 60

    # "2026-02-17" looks like a date, so Excel would normally convert it to a
    # date serial on assignment. Force text storage, then restore the
    # default "Normal" style so no extra number-format style sticks around.
    $ws.Range("B61").NumberFormat = "@"
    $ws.Range("B61").Value = "2026-02-17"
    $ws.Range("B61").Style = "Normal"

    # "12:52:49" is stored as plain text without needing any special
    # handling (it isn't auto-converted to a time serial).
    $ws.Range("C61").Value = "12:52:49"

    $ws.Range("D61").Value = "MarketMaking"
    $ws.Range("E61").Value = "DOWN"
    $ws.Range("F61").Value = 0.83
    $ws.Range("G61").Value = 0.85
    $ws.Range("H61").Value = "CLOSED"
    $ws.Range("I61").Value = 2.4096
    $ws.Range("J61").Value = 0.02
    $ws.Range("K61").Value = 100.1
    $ws.Range("L61").Value = 0
    $ws.Range("M61").Value = 0
    $ws.Range("N61").Value = 0.6
    $ws.Range("O61").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P61").Value = "early_exit"
    $ws.Range("Q61").Value = 0.11
}
